# Applies the "Penalty Reward System" forecast update:
#  - Shifts the Week_Start_Date column (B) forward by one week for each row
#    on the "Forecast Comparison" sheet
#  - Updates the MyForecast column (D) with new forecast values
#  - Refreshes the derived summary statistics on the "Summary" sheet
#
# NOTE: several of these values look like dates/numbers but must stay plain
# text (matching the source file's inline-string cells). A leading
# apostrophe is used - same as typing '2025-01-12 directly into Excel -
# so the values are stored as text instead of being auto-converted to a
# date/number serial.

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# New Week_Start_Date values (shifted one week later) for rows 2-17
$newDates = @(
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20",
    "2025-04-27"
)

# New MyForecast (column D) values for rows 2-17
$newForecast = @(48, 45, 41, 40, 45, 49, 52, 54, 55, 60, 40, 42, 56, 46, 41, 41)

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $i + 2
    $wsForecast.Cells.Item($row, 2).Value = "'" + $newDates[$i]
    $wsForecast.Cells.Item($row, 4).Value = $newForecast[$i]
}

# Update the Summary sheet with the refreshed statistics
$wsSummary.Range("B2").Value  = "2023-02-19 to 2025-01-05"
$wsSummary.Range("B3").Value  = "'0"
$wsSummary.Range("B4").Value  = "'87"
$wsSummary.Range("B5").Value  = "'44"
$wsSummary.Range("B7").Value  = "'18"
$wsSummary.Range("B8").Value  = "2888 units"
$wsSummary.Range("B9").Value  = "'756"
$wsSummary.Range("B10").Value = "'374"
$wsSummary.Range("B11").Value = "'174"
$wsSummary.Range("B12").Value = "'60"
$wsSummary.Range("B13").Value = "'2025-03-16"
$wsSummary.Range("B14").Value = "'40"
$wsSummary.Range("B15").Value = "'2025-02-02"
